# Fix "mont's pythons" -> "monty's pythons" typo (value not referenced by any cell,
# but corrected for data hygiene) and "x code" -> "x codes" are legacy/unused
# shared-string fixes; the visible effect on the grid comes from the new
# attendance columns added below for August 20 2016 and August 21 2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# New date headers for the two added day columns
$ws.Range("AC1").Value = "August 20 2016"
$ws.Range("AD1").Value = "August 21 2016"

# monty's pythons (row 2) clocked in/out on the new dates
$ws.Range("AC2").Value = "05:49 PM"
$ws.Range("AD2").Value = "07:08 PM"

# ruby perl (row 5) clocked in/out on the new dates
$ws.Range("AC5").Value = "05:51 PM"
$ws.Range("AD5").Value = "07:08 PM"

# sql database (row 8) clocked in/out on the new dates
$ws.Range("AC8").Value = "05:50 PM"
$ws.Range("AD8").Value = "07:10 PM"

# node .js (row 9) clocked in on the first new date only
$ws.Range("AC9").Value = "05:49 PM"

# a d (row 10) clocked in/out on the new dates
$ws.Range("AC10").Value = "05:50 PM"
$ws.Range("AD10").Value = "07:11 PM"
